$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 3 new columns at K:M for num_decode_1/2/3 (old K,L,M,N shift to N,O,P,Q)
$ws.Range("K1:M1").EntireColumn.Insert()

# 2. New header cells in the freshly inserted columns (order controls shared-string id order)
$ws.Range("K1").Value = "num_decode_1"
$ws.Range("L1").Value = "num_decode_2"
$ws.Range("M1").Value = "num_decode_3"

# 3. Fix F2 value (124 -> 128)
$ws.Range("F2").Value = 128

# 4. New data for row 2 in inserted columns
$ws.Range("K2").Value = 128
$ws.Range("L2").Value = 256
$ws.Range("M2").Value = "NAN"

# New NOTES header column (R)
$ws.Range("R1").Value = "NOTES"

# 5. Entire new row 3
$ws.Range("A3").Value = 5000
$ws.Range("B3").Value = 0.05
$ws.Range("C3").Value = 2048
$ws.Range("D3").Value = 1000
$ws.Range("E3").Value = 128
$ws.Range("F3").Value = 64
$ws.Range("G3").Value = 32
$ws.Range("H3").Value = 16
$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 8
$ws.Range("K3").Value = 32
$ws.Range("L3").Value = 64
$ws.Range("M3").Value = 128
$ws.Range("N3").Value = 8
$ws.Range("O3").Value = 0.3
$ws.Range("P3").Value = 0.0005
$ws.Range("R3").Value = "Increasing alpha helps training a lot"
$ws.Range("Q3").Value = "0.65 (sometimes 0.72)"

# 6. Apply center/center alignment across the whole used range (but R2 stays
# untouched/empty - only A1:Q3 plus R1 and R3 are populated cells).
# Doing this on a single cell first keeps the style table minimal (one new xf),
# then Copy + PasteSpecial(Formats) propagates that exact style everywhere
# without minting extra (orphan) style records.
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1").VerticalAlignment = -4108
$ws.Range("A1").Copy()
$ws.Range("A1:Q3").PasteSpecial(-4122)
$ws.Range("R1").PasteSpecial(-4122)
$ws.Range("R3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 7. Column widths for the newly introduced columns (approximate the
# bestFit widths the original authoring Excel session computed - this
# engine quantizes ColumnWidth to 1/6-character steps, so these are the
# closest achievable values to the authored widths of 14.7109375/
# 23.7109375/39.140625 characters respectively).
$ws.Range("K1:M1").ColumnWidth = 13.8333333333334
$ws.Range("Q1").ColumnWidth = 22.8333333333333
$ws.Range("R1").ColumnWidth = 38.3333333333334

# 8. Restore view state: scroll so column L is leftmost, select Q7.
$ws.Range("Q7").Select() | Out-Null
